# Convert the "Tamaño real del tráfico" column (D) from text values like
# "1,9 kB" / "52,7 kB" / "100 MB" into plain numeric values expressed in MB,
# shrink the now-unwrapped rows back down to their natural height, widen
# column A to fit the long file names, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D column: replace textual sizes with numeric MB values -----------------
$ws.Range("D2").Value = 0.0019
$ws.Range("D3").Value = 0.0527
$ws.Range("D4").Value = 0.3327
$ws.Range("D5").Value = 100
$ws.Range("D6").Value = 3.9
$ws.Range("D7").Value = 103
$ws.Range("D8").Value = 109
$ws.Range("D9").Value = 107
$ws.Range("D10").Value = 106
$ws.Range("D11").Value = 48

# --- Row heights: rows shrink now that column C/D text no longer wraps ------
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 28.5
$ws.Rows.Item(6).RowHeight = 28.5
$ws.Rows.Item(7).RowHeight = 28.5
$ws.Rows.Item(8).RowHeight = 28.5
$ws.Rows.Item(9).RowHeight = 28.5
$ws.Rows.Item(10).RowHeight = 28.5
$ws.Rows.Item(11).RowHeight = 28.5

# --- Column A: widen to fit the long pcap file names -------------------------
$ws.Columns.Item(1).ColumnWidth = 43.8

# --- Move the active selection ------------------------------------------------
$ws.Range("D12").Select() | Out-Null
